$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AprilRaw")

# Header row
$ws.Range("A1").Value2 = "Library"
$ws.Range("B1").Value2 = "Items owned by this library checked out at this library this month"
$ws.Range("C1").Value2 = "Items owned by other libraries checked out at this library this month"
$ws.Range("D1").Value2 = "Total circulation this month"

$ws.Range("A2").Value2 = "Atchison Public Library"
$ws.Range("B2").Value2 = 3847
$ws.Range("C2").Value2 = 1231
$ws.Range("D2").Value2 = 5078

$ws.Range("A3").Value2 = "Baldwin City Public Library"
$ws.Range("B3").Value2 = 2187
$ws.Range("C3").Value2 = 448
$ws.Range("D3").Value2 = 2635

$ws.Range("A4").Value2 = "Basehor Community Library"
$ws.Range("B4").Value2 = 8047
$ws.Range("C4").Value2 = 1362
$ws.Range("D4").Value2 = 9409

$ws.Range("A5").Value2 = "Bern Community Library"
$ws.Range("B5").Value2 = 145
$ws.Range("C5").Value2 = 73
$ws.Range("D5").Value2 = 218

$ws.Range("A6").Value2 = "Bonner Springs City Library"
$ws.Range("B6").Value2 = 4986
$ws.Range("C6").Value2 = 1367
$ws.Range("D6").Value2 = 6353

$ws.Range("A7").Value2 = "Burlingame Community Library"
$ws.Range("B7").Value2 = 554
$ws.Range("C7").Value2 = 282
$ws.Range("D7").Value2 = 836

$ws.Range("A8").Value2 = "Carbondale City Library"
$ws.Range("B8").Value2 = 577
$ws.Range("C8").Value2 = 134
$ws.Range("D8").Value2 = 711

$ws.Range("A9").Value2 = "Centralia Community Library"
$ws.Range("B9").Value2 = 321
$ws.Range("C9").Value2 = 75
$ws.Range("D9").Value2 = 396

$ws.Range("A10").Value2 = "Corning City Library"
$ws.Range("B10").Value2 = 27
$ws.Range("C10").Value2 = 6
$ws.Range("D10").Value2 = 33

$ws.Range("A11").Value2 = "Digital Content"

$ws.Range("A12").Value2 = "Doniphan County Library - Elwood"
$ws.Range("B12").Value2 = 104
$ws.Range("C12").Value2 = 13
$ws.Range("D12").Value2 = 117

$ws.Range("A13").Value2 = "Doniphan County Library - Highland"
$ws.Range("B13").Value2 = 254
$ws.Range("C13").Value2 = 91
$ws.Range("D13").Value2 = 345

$ws.Range("A14").Value2 = "Doniphan County Library - Troy"
$ws.Range("B14").Value2 = 549
$ws.Range("C14").Value2 = 129
$ws.Range("D14").Value2 = 678

$ws.Range("A15").Value2 = "Doniphan County Library - Wathena"
$ws.Range("B15").Value2 = 427
$ws.Range("C15").Value2 = 93
$ws.Range("D15").Value2 = 520

$ws.Range("A16").Value2 = "Effingham Community Library"
$ws.Range("B16").Value2 = 232
$ws.Range("C16").Value2 = 32
$ws.Range("D16").Value2 = 264

$ws.Range("A17").Value2 = "Eudora Community Library"
$ws.Range("B17").Value2 = 1677
$ws.Range("C17").Value2 = 612
$ws.Range("D17").Value2 = 2289

$ws.Range("A18").Value2 = "Everest, Barnes Reading Room"
$ws.Range("B18").Value2 = 118
$ws.Range("C18").Value2 = 172
$ws.Range("D18").Value2 = 290

$ws.Range("A19").Value2 = "Hiawatha, Morrill Public Library"
$ws.Range("B19").Value2 = 1595
$ws.Range("C19").Value2 = 635
$ws.Range("D19").Value2 = 2230

$ws.Range("A20").Value2 = "Highland Community College"
$ws.Range("B20").Value2 = 18
$ws.Range("C20").Value2 = 1
$ws.Range("D20").Value2 = 19

$ws.Range("A21").Value2 = "Holton, Beck-Bookman Library"
$ws.Range("B21").Value2 = 1601
$ws.Range("C21").Value2 = 506
$ws.Range("D21").Value2 = 2107

$ws.Range("A22").Value2 = "Horton Public Library"
$ws.Range("B22").Value2 = 103
$ws.Range("C22").Value2 = 13
$ws.Range("D22").Value2 = 116

$ws.Range("A23").Value2 = "Lansing Community Library"
$ws.Range("B23").Value2 = 1752
$ws.Range("C23").Value2 = 711
$ws.Range("D23").Value2 = 2463

$ws.Range("A24").Value2 = "Leavenworth Public Library"
$ws.Range("B24").Value2 = 8784
$ws.Range("C24").Value2 = 1897
$ws.Range("D24").Value2 = 10681

$ws.Range("A25").Value2 = "Linwood Community Library"
$ws.Range("B25").Value2 = 426
$ws.Range("C25").Value2 = 159
$ws.Range("D25").Value2 = 585

$ws.Range("A26").Value2 = "Louisburg Library"

$ws.Range("A27").Value2 = "Lyndon Carnegie Library"
$ws.Range("B27").Value2 = 358
$ws.Range("C27").Value2 = 270
$ws.Range("D27").Value2 = 628

$ws.Range("A28").Value2 = "McLouth Public Library"
$ws.Range("B28").Value2 = 143
$ws.Range("C28").Value2 = 63
$ws.Range("D28").Value2 = 206

$ws.Range("A29").Value2 = "Meriden-Ozawkie Public Library"
$ws.Range("B29").Value2 = 1041
$ws.Range("C29").Value2 = 457
$ws.Range("D29").Value2 = 1498

$ws.Range("A30").Value2 = "Northeast Kansas Library System"
$ws.Range("B30").Value2 = 26
$ws.Range("C30").Value2 = 21
$ws.Range("D30").Value2 = 47

$ws.Range("A31").Value2 = "Nortonville Public Library"
$ws.Range("B31").Value2 = 294
$ws.Range("C31").Value2 = 111
$ws.Range("D31").Value2 = 405

$ws.Range("A32").Value2 = "Osage City Library"
$ws.Range("B32").Value2 = 1582
$ws.Range("C32").Value2 = 404
$ws.Range("D32").Value2 = 1986

$ws.Range("A33").Value2 = "Osawatomie Public Library"
$ws.Range("B33").Value2 = 774
$ws.Range("C33").Value2 = 388
$ws.Range("D33").Value2 = 1162

$ws.Range("A34").Value2 = "Oskaloosa Public Library"
$ws.Range("B34").Value2 = 454
$ws.Range("C34").Value2 = 222
$ws.Range("D34").Value2 = 676

$ws.Range("A35").Value2 = "Ottawa Library"
$ws.Range("B35").Value2 = 5844
$ws.Range("C35").Value2 = 889
$ws.Range("D35").Value2 = 6733

$ws.Range("A36").Value2 = "Overbrook Public Library"
$ws.Range("B36").Value2 = 871
$ws.Range("C36").Value2 = 225
$ws.Range("D36").Value2 = 1096

$ws.Range("A37").Value2 = "Paola Free Library"
$ws.Range("B37").Value2 = 2670
$ws.Range("C37").Value2 = 477
$ws.Range("D37").Value2 = 3147

$ws.Range("A38").Value2 = "Perry-Lecompton Community Library"
$ws.Range("B38").Value2 = 63
$ws.Range("C38").Value2 = 33
$ws.Range("D38").Value2 = 96

$ws.Range("A39").Value2 = "Pomona Community Library"
$ws.Range("B39").Value2 = 119
$ws.Range("C39").Value2 = 130
$ws.Range("D39").Value2 = 249

$ws.Range("A40").Value2 = "Prairie Hills Schools - Axtell Public School"
$ws.Range("B40").Value2 = 449
$ws.Range("C40").Value2 = 2
$ws.Range("D40").Value2 = 451

$ws.Range("A41").Value2 = "Prairie Hills Schools - Sabetha Elementary School"
$ws.Range("B41").Value2 = 1488
$ws.Range("C41").Value2 = 92
$ws.Range("D41").Value2 = 1580

$ws.Range("A42").Value2 = "Prairie Hills Schools - Sabetha High School"
$ws.Range("B42").Value2 = 24
$ws.Range("C42").Value2 = 5
$ws.Range("D42").Value2 = 29

$ws.Range("A43").Value2 = "Prairie Hills Schools - Sabetha Middle School"
$ws.Range("B43").Value2 = 118
$ws.Range("C43").Value2 = 5
$ws.Range("D43").Value2 = 123

$ws.Range("A44").Value2 = "Prairie Hills Schools - Wetmore Academic Center"
$ws.Range("B44").Value2 = 194
$ws.Range("D44").Value2 = 194

$ws.Range("A45").Value2 = "Richmond Public Library"
$ws.Range("B45").Value2 = 396
$ws.Range("C45").Value2 = 74
$ws.Range("D45").Value2 = 470

$ws.Range("A46").Value2 = "Rossville Community Library"
$ws.Range("B46").Value2 = 1176
$ws.Range("C46").Value2 = 368
$ws.Range("D46").Value2 = 1544

$ws.Range("A47").Value2 = "Sabetha, Mary Cotton Library"
$ws.Range("B47").Value2 = 3059
$ws.Range("C47").Value2 = 915
$ws.Range("D47").Value2 = 3974

$ws.Range("A48").Value2 = "Seneca Free Library"
$ws.Range("B48").Value2 = 1426
$ws.Range("C48").Value2 = 150
$ws.Range("D48").Value2 = 1576

$ws.Range("A49").Value2 = "Silver Lake Library"
$ws.Range("B49").Value2 = 1107
$ws.Range("C49").Value2 = 349
$ws.Range("D49").Value2 = 1456

$ws.Range("A50").Value2 = "Tonganoxie Public Library"
$ws.Range("B50").Value2 = 2822
$ws.Range("C50").Value2 = 845
$ws.Range("D50").Value2 = 3667

$ws.Range("A51").Value2 = "Valley Falls, Delaware Township Library"
$ws.Range("B51").Value2 = 409
$ws.Range("C51").Value2 = 300
$ws.Range("D51").Value2 = 709

$ws.Range("A52").Value2 = "Wellsville City Library"
$ws.Range("B52").Value2 = 1279
$ws.Range("C52").Value2 = 514
$ws.Range("D52").Value2 = 1793

$ws.Range("A53").Value2 = "Wetmore Public Library"
$ws.Range("B53").Value2 = 111
$ws.Range("C53").Value2 = 97
$ws.Range("D53").Value2 = 208

$ws.Range("A54").Value2 = "Williamsburg Community Library"
$ws.Range("B54").Value2 = 289
$ws.Range("C54").Value2 = 23
$ws.Range("D54").Value2 = 312

$ws.Range("A55").Value2 = "Winchester Public Library"
$ws.Range("B55").Value2 = 239
$ws.Range("C55").Value2 = 151
$ws.Range("D55").Value2 = 390
